$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.952.51"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.916.95"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "3.401.93"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "60.927.06"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.68"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "2.917.11"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "430.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.679"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D34").Value = "0.0₃0854"
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.02"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.54"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.92"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.20"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "2.701.09"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "131.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("E51").Value = "  +2.20%  "
